# "Update countries & provincias Spain"
#
# Refresh the COVID dashboard snapshot on sheet "Pais":
#   - bump the "last updated" timestamp banner in A1
#   - refresh per-country counters (Casos totales / Nuevos casos /
#     Casos activos / Recuperados / Casos criticos / Muertes hoy / Muertes)
#     for the countries whose numbers moved in the new data pull. The
#     table is kept sorted by "Casos totales" (column B) descending, so a
#     few neighbouring countries swap places:
#       - Kirguistan now outranks Paraguay (rows 67/68)
#       - Malasia jumps above Noruega/Zambia/Senegal/Albania (rows 94-98)
#       - Gibraltar now outranks Islas Feroe (rows 180/181)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Datos actualizados a 11 de Octubre de 2020 a las 11:20" -> "...12:37"
$ws.Range("A1").Value = "Datos actualizados a 11 de Octubre de 2020 a las 12:37"

# row -> { Pais, Casos totales, Nuevos casos, Casos activos, Recuperados,
#          Casos criticos, Muertes hoy, Muertes }
$countryData = @{
    19  = @{ Pais = "Banglades";   CasosTotales = 378266; NuevosCasos = 1193; CasosActivos = 292860; Recuperados = 79882; CasosCriticos = 0; MuertesHoy = 24; Muertes = 5524 }
    64  = @{ Pais = "Austria";     CasosTotales = 55319;  NuevosCasos = 896;  CasosActivos = 43448;  Recuperados = 11020; CasosCriticos = 0; MuertesHoy = 0;  Muertes = 851 }
    67  = @{ Pais = "Kirguistan";  CasosTotales = 49230;  NuevosCasos = 306;  CasosActivos = 44227;  Recuperados = 3918;  CasosCriticos = 0; MuertesHoy = 3;  Muertes = 1085 }
    68  = @{ Pais = "Paraguay";    CasosTotales = 48978;  NuevosCasos = 0;    CasosActivos = 31351;  Recuperados = 16562; CasosCriticos = 0; MuertesHoy = 0;  Muertes = 1065 }
    94  = @{ Pais = "Malasia";     CasosTotales = 15657;  NuevosCasos = 561;  CasosActivos = 10913;  Recuperados = 4587;  CasosCriticos = 0; MuertesHoy = 2;  Muertes = 157 }
    95  = @{ Pais = "Noruega";     CasosTotales = 15466;  NuevosCasos = 0;    CasosActivos = 11863;  Recuperados = 3328;  CasosCriticos = 0; MuertesHoy = 0;  Muertes = 275 }
    96  = @{ Pais = "Zambia";      CasosTotales = 15415;  NuevosCasos = 0;    CasosActivos = 14541;  Recuperados = 537;   CasosCriticos = 0; MuertesHoy = 0;  Muertes = 337 }
    97  = @{ Pais = "Senegal";     CasosTotales = 15244;  NuevosCasos = 0;    CasosActivos = 13198;  Recuperados = 1732;  CasosCriticos = 0; MuertesHoy = 0;  Muertes = 314 }
    98  = @{ Pais = "Albania";     CasosTotales = 15231;  NuevosCasos = 0;    CasosActivos = 9406;   Recuperados = 5409;  CasosCriticos = 0; MuertesHoy = 0;  Muertes = 416 }
    121 = @{ Pais = "Lituania";    CasosTotales = 6122;   NuevosCasos = 160;  CasosActivos = 2777;   Recuperados = 3242;  CasosCriticos = 0; MuertesHoy = 0;  Muertes = 103 }
    127 = @{ Pais = "Hong Kong";   CasosTotales = 5183;   NuevosCasos = 7;    CasosActivos = 4919;   Recuperados = 159;   CasosCriticos = 0; MuertesHoy = 0;  Muertes = 105 }
    142 = @{ Pais = "Malta";       CasosTotales = 3776;   NuevosCasos = 95;   CasosActivos = 2967;   Recuperados = 768;   CasosCriticos = 0; MuertesHoy = 0;  Muertes = 41 }
    180 = @{ Pais = "Gibraltar";   CasosTotales = 485;    NuevosCasos = 9;    CasosActivos = 413;    Recuperados = 72;    CasosCriticos = 0; MuertesHoy = 0;  Muertes = 0 }
    181 = @{ Pais = "Islas Feroe"; CasosTotales = 477;    NuevosCasos = 0;    CasosActivos = 461;    Recuperados = 16;    CasosCriticos = 0; MuertesHoy = 0;  Muertes = 0 }
}

foreach ($row in $countryData.Keys) {
    $d = $countryData[$row]
    $ws.Cells.Item($row, 1).Value = $d.Pais
    $ws.Cells.Item($row, 2).Value = $d.CasosTotales
    $ws.Cells.Item($row, 3).Value = $d.NuevosCasos
    $ws.Cells.Item($row, 4).Value = $d.CasosActivos
    $ws.Cells.Item($row, 5).Value = $d.Recuperados
    $ws.Cells.Item($row, 6).Value = $d.CasosCriticos
    $ws.Cells.Item($row, 7).Value = $d.MuertesHoy
    $ws.Cells.Item($row, 8).Value = $d.Muertes
}
